# The "Product In Market" service row (row 8) had its Resource Description
# updated from "Supermercados" to "Supermercados, Produtos".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "Supermercados, Produtos"

$wb.Save()
